# Harmonize terminology in the header row of the dims_structure sheet
# (use type -> function type, construction type -> structure type, etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "function_type"
$ws.Range("E1").Value = "structure_type_subtype"
$ws.Range("F1").Value = "energy_efficiency"
$ws.Range("H1").Value = "function"
$ws.Range("I1").Value = "structure"

# Update the saved selection/active cell for the sheet
$ws.Range("J1").Select()
